# Append the new "2021年" row (row 11) to Sheet1, mirroring the layout of
# the existing yearly rows (row 10 = "2020年").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A11").Value = "2021年"
$ws.Range("B11").Value = 890
$ws.Range("C11").Value = 259
$ws.Range("D11").Value = 60

# E11 mirrors E10: an explicit empty text cell (not a blank/number cell).
# Entering a lone apostrophe forces Excel/COM to store it as text with an
# empty value; the quote-prefix formatting that trick adds is then
# stripped by restoring the cell's style to match its neighbours.
$ws.Range("E11").Value = "'"
$ws.Range("E11").Style = $ws.Range("F11").Style

$ws.Range("F11").Value = 816
$ws.Range("G11").Value = 1271
$ws.Range("H11").Value = 67
$ws.Range("I11").Value = 540
$ws.Range("J11").Value = 291
$ws.Range("K11").Value = 25180
$ws.Range("L11").Value = 23
$ws.Range("M11").Value = 171
$ws.Range("N11").Value = 41
$ws.Range("O11").Value = 68
$ws.Range("P11").Value = 586
$ws.Range("Q11").Value = 262
$ws.Range("R11").Value = 66
$ws.Range("S11").Value = 276
$ws.Range("T11").Value = 1967
$ws.Range("U11").Value = 889
$ws.Range("V11").Value = 81
$ws.Range("W11").Value = 990
$ws.Range("X11").Value = 874
$ws.Range("Y11").Value = 6192
$ws.Range("Z11").Value = 719
$ws.Range("AA11").Value = 16
$ws.Range("AB11").Value = 268
$ws.Range("AC11").Value = 101
$ws.Range("AD11").Value = 152
$ws.Range("AE11").Value = 226
$ws.Range("AF11").Value = 883
$ws.Range("AG11").Value = 833
$ws.Range("AH11").Value = 97
$ws.Range("AI11").Value = 329
$ws.Range("AJ11").Value = 154
$ws.Range("AK11").Value = 642
$ws.Range("AL11").Value = 630
$ws.Range("AM11").Value = 2234
$ws.Range("AN11").Value = 413
$ws.Range("AO11").Value = 380
$ws.Range("AP11").Value = 293
$ws.Range("AQ11").Value = 129

# Match the style of the equivalent cell in the prior year's row (A10 is
# bold/bordered; the data columns use the default style). Copy/PasteSpecial
# (formats only) reuses the existing style record instead of synthesizing a
# new font, matching how the rest of column A is formatted.
$ws.Range("A10").Copy()
$ws.Range("A11").PasteSpecial(-4122)
$excel.CutCopyMode = $false
